# Apply alpha-diversity stat corrections to the single data table.
#
# Many of the affected values (e.g. "1.85", "0.068", "0.069", "0.479")
# occur more than once in the document, so every replacement is scoped to
# the exact Table.Cell(row, column) that the XML diff shows changing,
# rather than relying on a document-wide Find/Replace.
#
# The "Estimate" / "Standard Error" columns store the visible number as one
# run immediately followed by a separate superscript exponent run (e.g.
# "8.1" + NBSP + "x" + NBSP + "10" in one run, "2" in the next run,
# together rendered as "8.1 x 10^2"). For those we use a cell-scoped
# Find/Replace (wdReplaceOne, *not* wdReplaceAll -- ReplaceAll silently
# escapes the cell scope in this runtime) so only the matched run's text
# changes and the superscript run is left untouched, mirroring the XML
# diff exactly.
#
# The "T Value" / "p" columns hold the whole cell value in a single run,
# so those are updated by assigning Range.Text directly, which is both
# simpler and avoids a Find engine quirk that can silently no-op a second
# identical Find/Replace pair run back-to-back.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nbsp = [char]0x00A0

function Set-CellText($row, $col, $newVal) {
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $newVal
}

function Set-CellRun($row, $col, $oldVal, $newVal) {
    $cell = $t.Cell($row, $col)
    $range = $cell.Range
    $found = $range.Find.Execute($oldVal, $true, $false, $false, $false, $false, $true, 0, $false, $newVal, 1)
    if (-not $found) {
        Write-Host "WARNING: replacement not found for row=$row col=$col old=[$oldVal] new=[$newVal]"
    }
}

# Row 14 (Diversity (Shannon H) / Intercept)
$old = "{0}{1}x{1}10" -f "8.1", $nbsp
$new = "{0}{1}x{1}10" -f "7.8", $nbsp
Set-CellRun 14 4 $old $new
Set-CellText 14 5 "1.89"
Set-CellText 14 6 "0.063"

# Row 15 (log(Size Class))
$old = "{0}{1}x{1}10" -f "4.9", $nbsp
$new = "{0}{1}x{1}10" -f "4.7", $nbsp
Set-CellRun 15 4 $old $new
Set-CellText 15 5 "4.21"

# Row 16 (log(Size Class)^2)
$old = "{0}{1}x{1}10" -f "9.1", $nbsp
$new = "{0}{1}x{1}10" -f "8.8", $nbsp
Set-CellRun 16 4 $old $new
Set-CellText 16 5 "-4.12"

# Row 17 (Latitude)
$old = "{0}{1}x{1}10" -f "-7.8", $nbsp
$new = "{0}{1}x{1}10" -f "-7.7", $nbsp
Set-CellRun 17 3 $old $new
$old = "{0}{1}x{1}10" -f "4.2", $nbsp
$new = "{0}{1}x{1}10" -f "4.1", $nbsp
Set-CellRun 17 4 $old $new
Set-CellText 17 5 "-1.88"
Set-CellText 17 6 "0.064"

# Row 18 (Latitude^2)
$old = "{0}{1}x{1}10" -f "5.5", $nbsp
$new = "{0}{1}x{1}10" -f "5.3", $nbsp
Set-CellRun 18 4 $old $new
Set-CellText 18 5 "1.88"
Set-CellText 18 6 "0.064"

# Row 19 (Depth)
Set-CellText 19 5 "0.82"
Set-CellText 19 6 "0.413"

# Row 20 (Evenness (Pielou J) / Intercept)
Set-CellText 20 6 "0.477"

# Row 23 (Latitude, under Evenness)
Set-CellText 23 6 "0.478"
